# Update the "想去人数" (number of people interested) counts for several
# conghui/events, on both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 47
$ws1.Range("F3").Value = 327
$ws1.Range("F5").Value = 3117
$ws1.Range("F6").Value = 2068
$ws1.Range("F7").Value = 397
$ws1.Range("F9").Value = 1163
$ws1.Range("F11").Value = 953
$ws1.Range("F12").Value = 81

# --- Sheet "全部类型" (all types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 47
$ws4.Range("F3").Value = 327
$ws4.Range("F5").Value = 3117
$ws4.Range("F6").Value = 2068
$ws4.Range("F7").Value = 397
$ws4.Range("F10").Value = 1163
$ws4.Range("F12").Value = 953
$ws4.Range("F13").Value = 81
